$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Target data for rows 2-7 (columns A..T)
# Column layout:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("ECs",  "Tnf", "Tnfrsf1a", "ECs",  3, 1,                  459.5553626666667, 1378.666088, 0.997342491800172,  0.997342491800172,  3, 1, 73.91316300000001,  221.739489, 0.6096331558809399, 0.6096331558809398, 33967.19042830545, 305704.713854749,  0.6080130507702993,  0.6080130507702992),
    @("ECs",  "Tnf", "Tnfrsf1a", "FAPs", 3, 1,                  459.5553626666667, 1378.666088, 0.997342491800172,  0.997342491800172,  3, 1, 31.40056566666667,  94.201697,  0.2589907557307037, 0.2589907557307037, 14430.29834288348, 129872.6850859513, 0.2583024856736697,  0.2583024856736697),
    @("ECs",  "Tnf", "Tnfrsf1a", "sCs",  3, 1,                  459.5553626666667, 1378.666088, 0.997342491800172,  0.997342491800172,  3, 1, 15.928304,          47.784912,  0.1313760883883564, 0.1313760883883564, 7319.937521384917, 65879.43769246426, 0.131026955356203,   0.131026955356203),
    @("FAPs", "Tnf", "Tnfrsf1a", "ECs",  2, 0.6666666666666666, 1.224526333333333, 3.673579,    0.002657508199827995,0.002657508199827995,3, 1, 73.91316300000001,  221.739489, 0.6096331558809399, 0.6096331558809398, 90.50861447345902, 814.5775302611311, 0.001620105110640616,0.001620105110640616),
    @("FAPs", "Tnf", "Tnfrsf1a", "FAPs", 2, 0.6666666666666666, 1.224526333333333, 3.673579,    0.002657508199827995,0.002657508199827995,3, 1, 31.40056566666667,  94.201697,  0.2589907557307037, 0.2589907557307037, 38.45081954039589, 346.057375863563,  0.0006882700570339943,0.0006882700570339943),
    @("FAPs", "Tnf", "Tnfrsf1a", "sCs",  2, 0.6666666666666666, 1.224526333333333, 3.673579,    0.002657508199827995,0.002657508199827995,3, 1, 15.928304,          47.784912,  0.1313760883883564, 0.1313760883883564, 19.50462769333867, 175.541649240048,  0.0003491330321533847,0.0003491330321533847)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
